$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-sort the "Periodo Mora" column (E16:E20) into ascending order:
# previously 1808,1807,1806,1805,1804 -> now 1804,1805,1806,1807,1808
$ws.Range("E16").Value = "1804"
$ws.Range("E17").Value = "1805"
$ws.Range("E18").Value = "1806"
$ws.Range("E19").Value = "1807"
$ws.Range("E20").Value = "1808"
